$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: replace the (unprefixed) indicator-name strings with the
#     country-name strings (Kyrgyz / Russian / English) ---
$ws.Range("A5").Value = "Кыргыз Республикасы"
$ws.Range("B5").Value = "Кыргызская Республика"
$ws.Range("C5").Value = "Kyrgyz Republic "

# --- New column R: year 2023 header + its data value ---
$ws.Range("R4").Value = 2023
$ws.Range("R5").Value = 53.5

# Copy the formatting of the neighbouring (Q) column cells onto the new
# R cells so they pick up the same borders / number format / alignment.
$ws.Range("Q4:Q5").Copy()
$ws.Range("R4:R5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 5 is shorter now (single-line country names instead of the long
# wrapped indicator description), shrink its height to match.
$ws.Rows.Item(5).RowHeight = 21

# --- Columns A:C are unified to one common, narrower width ---
$ws.Columns("A:C").ColumnWidth = 35

# Reset the selection back to the top-left cell.
$ws.Range("A1").Select()
